# Added one input and one output (opto+relay).
# Use for Start/Selection and Main Heater Control.
#
# This adds a new "SW var name" column (J) to the pinout table, populating
# the variable names used for the new Usr_Btn_1/2/3 signals, resizes the
# embedded Visio drawing placeholder slightly, and updates the active
# selection in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J : "SW var name" -----------------------------------------
# Shared strings must end up in this order: Usr_Btn_3, Usr_Btn_2, Usr_Btn_1,
# SW var name -- so write the cells in that same order.
$ws.Range("J9").Value = "Usr_Btn_3"
$ws.Range("J6").Value = "Usr_Btn_2"
$ws.Range("J5").Value = "Usr_Btn_1"
$ws.Range("J3").Value = "SW var name"

# Match the header formatting used by the rest of row 3 (bold + wrap, style
# of cell I3).
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Size column J to fit its new contents.
$ws.Columns("J:J").AutoFit()

# --- Resize the embedded Visio object placeholder --------------------------
# The shape's right edge (xdr:to colOff in column L) moves from 485775 EMU
# to 438150 EMU (5 px to the left). The exact width below was derived from
# the shape's fixed absolute left offset (2189051 EMU) and the absolute EMU
# position where column L begins after the column J resize above
# (8899227 EMU), so that left + width lands exactly on the target offset.
$shp = $ws.Shapes.Item(1)
$shp.Width = 562.8603149606299

# --- Update the active selection / scrolled view ---------------------------
$ws.Range("J6").Select() | Out-Null
